# Update "想去人数" (F column) counts across sheets to reflect refreshed
# output data (regenerated gh-pages data at commit 802b57d).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Changes = @{
    2  = 979
    4  = 846
    6  = 6918
    9  = 266
    11 = 521
    12 = 8
    14 = 369
    16 = 2475
    17 = 117
    18 = 184
    19 = 734
    21 = 412
    22 = 82
    24 = 38
    25 = 100
    26 = 44
    27 = 150
    30 = 332
    31 = 15
}
foreach ($row in $sheet1Changes.Keys) {
    $ws1.Range("F$row").Value = $sheet1Changes[$row]
}

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 186

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Changes = @{
    2  = 186
    3  = 979
    6  = 846
    10 = 6918
    13 = 266
    15 = 521
    16 = 8
    18 = 369
    21 = 2475
    22 = 117
    24 = 184
    25 = 734
    28 = 412
    29 = 82
    31 = 38
    32 = 100
    33 = 44
    34 = 150
    37 = 332
    38 = 15
}
foreach ($row in $sheet4Changes.Keys) {
    $ws4.Range("F$row").Value = $sheet4Changes[$row]
}

$wb.Save()
